$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = 12
$ws.Range("J2").Value = 0.0008333333333333334
$ws.Range("K2").Value = 5691
$ws.Range("L2").Value = 0.011382
